# "add drop name column" - rename the generic "建筑N" (Building N) slot
# names to "附楼N" (Annex N) in the Slot table's Name column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "附楼1"
$ws.Range("B6").Value = "附楼2"
$ws.Range("B11").Value = "附楼3"

# Match the author's final cursor position recorded in the saved sheet view.
$ws.Range("B8").Select()
